$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Branch string in A3
$ws.Range("A3").Value = "CIMS.CAN.AB.Residential.Dwellings.Lighting"

# Replace formulas in N3:W3 with static values, and update M3
$ws.Range("M3").Value = 1
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 0.5
$ws.Range("Q3").Value = 0.0001
$ws.Range("R3").Value = 0.0001
$ws.Range("S3").Value = 0.0001
$ws.Range("T3").Value = 0.0001
$ws.Range("U3").Value = 0.0001
$ws.Range("V3").Value = 0.0001
$ws.Range("W3").Value = 0.0001

# Update selected range in the sheet view
$ws.Range("A1:X4").Select()
